# Daily attendance processing - 2025-11-24 07:48:46
# Applies the attendance-report refresh to the "Session Analysis Results" sheet:
#  - reorders the "Recorded By" name lists on several rows
#  - updates recorded-session counters / coverage percentages
#  - flips the MICROBIOLOGY C1 session #1 (row 12) from "Not Recorded" to "Recorded"
#    now that it has been recorded by yassmina.fattoh@med.asu.edu.eg

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 (ANATOMY, session 1): reorder "Recorded By" list ----
$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, System"

# ---- Row 3 (ANATOMY, session 2): reorder "Recorded By" list ----
$ws.Range("G3").Value = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# ---- Row 4 (ANATOMY, session 3): reorder "Recorded By" list ----
$ws.Range("G4").Value = "hend_mahmoud@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# ---- Row 6: Recorded Sessions metric 11 -> 12 ----
$ws.Range("L6").Value = 12

# ---- Row 7 (BIOCHEMISTRY LAB/CBL, session 1): reorder "Recorded By" list ----
$ws.Range("G7").Value = "menna-alah.mohamed@asu.edu.eg, AbeerRagheb@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg"

# ---- Row 7: Missing Sessions metric 2 -> 1 ----
$ws.Range("L7").Value = 1

# ---- Row 9 (HISTOLOGY, session 1): reorder "Recorded By" list ----
$ws.Range("G9").Value = "Shimaa.ashraf@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"

# ---- Row 9: Coverage % 37.9% -> 41.4% (keep as literal text, preserving original cell style) ----
$ws.Range("ZZ1").Formula = '="41.4%"'
$ws.Range("ZZ1").Copy()
$ws.Range("L9").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()
$excel.CutCopyMode = $false

# ---- Row 10: Average Attendance % 24.9% -> 23.6% (keep as literal text, preserving original cell style) ----
$ws.Range("ZZ1").Formula = '="23.6%"'
$ws.Range("ZZ1").Copy()
$ws.Range("L10").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()
$excel.CutCopyMode = $false

# ---- Row 12 (MICROBIOLOGY C1, session 1): now recorded ----
$ws.Range("G12").Value = "yassmina.fattoh@med.asu.edu.eg"
$ws.Range("H12").Value = "24/251"
$ws.Range("I12").Value = "Recorded"
# Re-color the row from the "Not Recorded" look to the "Recorded" look (copy formats from a Recorded row)
$ws.Range("A2:I2").Copy()
$ws.Range("A12:I12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Row 15 (PARASITOLOGY class-statistics summary row): refresh rolled-up metrics ----
$ws.Range("O15").Value = 12
$ws.Range("P15").Value = 1

$ws.Range("ZZ1").Formula = '="41.4%"'
$ws.Range("ZZ1").Copy()
$ws.Range("R15").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()
$excel.CutCopyMode = $false

$ws.Range("ZZ1").Formula = '="23.6%"'
$ws.Range("ZZ1").Copy()
$ws.Range("S15").PasteSpecial(-4163)
$ws.Range("ZZ1").Clear()
$excel.CutCopyMode = $false
